$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 9;  Year = 2003; Japanese = "コール オブ クトゥルフ d20"; English = "Call of Cthulhu d20"; Publisher = "Shinkigensha"; Image = "d20_rulebook.jpg"; ProductType = "rulebook" },
    @{ Row = 10; Year = 2004; Japanese = "H.P.ラヴクラフト アーカム"; English = "H.P. Lovecraft's Arkham"; Publisher = "Shinkigensha"; Image = "arkham.jpg"; ProductType = "supplement" },
    @{ Row = 11; Year = 2005; Japanese = "クトゥルフ・ダークエイジ"; English = "Cthulhu Dark Ages"; Publisher = "Shinkigensha"; Image = "dark_ages.jpg"; ProductType = "supplement" },
    @{ Row = 12; Year = 2006; Japanese = "クトゥルフ神話TRPG 比叡山炎上"; English = "Hieizan Flame"; Publisher = "Enterbrain"; Image = "hieizan_flame.jpg"; ProductType = "supplement" },
    @{ Row = 13; Year = 2007; Japanese = "クトゥルフ神話TRPGシナリオ集 七つの怪談"; English = "Seven Phantom Stories"; Publisher = "Shinkigensha"; Image = "seven_phantom_stories.jpg"; ProductType = "scenario" },
    @{ Row = 14; Year = 2008; Japanese = "クトゥルフ神話TRPG マレウス・モンストロルム"; English = "Malleus Monstrorum"; Publisher = "Enterbrain"; Image = "malleus_monstrorum.jpg"; ProductType = "supplement" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Year
    $ws.Cells.Item($row, 2).Value = $r.Japanese
    $ws.Cells.Item($row, 3).Value = $r.English
    $ws.Cells.Item($row, 4).Value = $r.Publisher
    $ws.Cells.Item($row, 5).Value = $r.Image
    $ws.Cells.Item($row, 6).Value = $r.ProductType
}

$ws.Range("F15").Select()
